# Insert 3 new data rows at the top of the "261" block, pushing the
# existing rows 261-321 down to 264-324 (dimension grows from T321 to T324).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("261:263").Insert()

# Shared/static columns for every data row in this sheet.
$mercado = "Agrícola del Norte S.A. de Arica"
$region  = "Arica y Parinacota"
$codreg  = 15
$tipo    = "Fruta"
$prodId  = 100108
$prod    = "Tropicales y subtropicales"
$catId   = 100108006
$cat     = "Plátano"
$unidad  = "$/caja 20 kilos"
$kgUnid  = 20

# Row 261 - new record
$ws.Cells.Item(261, 1).Value = 1
$ws.Cells.Item(261, 2).Value = $mercado
$ws.Cells.Item(261, 3).Value = $region
$ws.Cells.Item(261, 4).Value = 44889
$ws.Cells.Item(261, 5).Value = $codreg
$ws.Cells.Item(261, 6).Value = $tipo
$ws.Cells.Item(261, 7).Value = $prodId
$ws.Cells.Item(261, 8).Value = $prod
$ws.Cells.Item(261, 9).Value = $catId
$ws.Cells.Item(261, 10).Value = $cat
$ws.Cells.Item(261, 11).Value = "Barraganete"
$ws.Cells.Item(261, 12).Value = "Primera"
$ws.Cells.Item(261, 13).Value = 950
$ws.Cells.Item(261, 14).Value = 30000
$ws.Cells.Item(261, 15).Value = 31000
$ws.Cells.Item(261, 16).Value = 30421
$ws.Cells.Item(261, 17).Value = $unidad
$ws.Cells.Item(261, 18).Value = "Ecuador"
$ws.Cells.Item(261, 19).Value = 1521
$ws.Cells.Item(261, 20).Value = $kgUnid

# Row 262 - new record
$ws.Cells.Item(262, 1).Value = 1
$ws.Cells.Item(262, 2).Value = $mercado
$ws.Cells.Item(262, 3).Value = $region
$ws.Cells.Item(262, 4).Value = 44889
$ws.Cells.Item(262, 5).Value = $codreg
$ws.Cells.Item(262, 6).Value = $tipo
$ws.Cells.Item(262, 7).Value = $prodId
$ws.Cells.Item(262, 8).Value = $prod
$ws.Cells.Item(262, 9).Value = $catId
$ws.Cells.Item(262, 10).Value = $cat
$ws.Cells.Item(262, 11).Value = "Sin especificar"
$ws.Cells.Item(262, 12).Value = "Pintón"
$ws.Cells.Item(262, 13).Value = 500
$ws.Cells.Item(262, 14).Value = 30000
$ws.Cells.Item(262, 15).Value = 31000
$ws.Cells.Item(262, 16).Value = 30700
$ws.Cells.Item(262, 17).Value = $unidad
$ws.Cells.Item(262, 18).Value = "Ecuador"
$ws.Cells.Item(262, 19).Value = 1535
$ws.Cells.Item(262, 20).Value = $kgUnid

# Row 263 - new record
$ws.Cells.Item(263, 1).Value = 1
$ws.Cells.Item(263, 2).Value = $mercado
$ws.Cells.Item(263, 3).Value = $region
$ws.Cells.Item(263, 4).Value = 44889
$ws.Cells.Item(263, 5).Value = $codreg
$ws.Cells.Item(263, 6).Value = $tipo
$ws.Cells.Item(263, 7).Value = $prodId
$ws.Cells.Item(263, 8).Value = $prod
$ws.Cells.Item(263, 9).Value = $catId
$ws.Cells.Item(263, 10).Value = $cat
$ws.Cells.Item(263, 11).Value = "Barraganete"
$ws.Cells.Item(263, 12).Value = "Verde"
$ws.Cells.Item(263, 13).Value = 340
$ws.Cells.Item(263, 14).Value = 30000
$ws.Cells.Item(263, 15).Value = 31000
$ws.Cells.Item(263, 16).Value = 30441
$ws.Cells.Item(263, 17).Value = $unidad
$ws.Cells.Item(263, 18).Value = "Ecuador"
$ws.Cells.Item(263, 19).Value = 1522
$ws.Cells.Item(263, 20).Value = $kgUnid
